# Add a new "onFormulaChanged" snippet entry to the Snippets table.
# This inserts two new rows (for the "registerFormulaChangeHandler" and
# "formulaChangeHandler" methods) right above the existing row that used
# to be row 218 ("Worksheet" / "showOutlineLevels" entry), and keeps the
# table/autofilter/dimension ranges in sync with the newly-added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows by inserting before the current row 218,
# which pushes all following rows (and the table's data) down by two rows.
$ws.Rows.Item(218).Insert()
$ws.Rows.Item(218).Insert()

# New row 218: Worksheet.onFormulaChanged -> registerFormulaChangeHandler
$ws.Cells.Item(218, 1).Value = "Worksheet"
$ws.Cells.Item(218, 2).Value = "onFormulaChanged"
$ws.Cells.Item(218, 4).Value = "excel-events-formula-changed"
$ws.Cells.Item(218, 5).Value = "registerFormulaChangeHandler"

# New row 219: Worksheet.onFormulaChanged -> formulaChangeHandler
$ws.Cells.Item(219, 1).Value = "Worksheet"
$ws.Cells.Item(219, 2).Value = "onFormulaChanged"
$ws.Cells.Item(219, 4).Value = "excel-events-formula-changed"
$ws.Cells.Item(219, 5).Value = "formulaChangeHandler"

# Grow the "Snippets" table/autofilter to include the two new rows so the
# table range covers A1:E248 instead of A1:E246.
$lo = $ws.ListObjects.Item(1)
$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1 + 2
$lo.Resize($ws.Range("A1:E" + $lastRow))

# Leave the selection on the newly added cell, matching where the author
# ended up after adding the snippet.
$ws.Range("C219").Select()
